$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.011.74"
$ws.Range("E2").Value = "  -0.03%  "

# Row 3
$ws.Range("D3").Value = "1.626.28"
$ws.Range("E3").Value = "  -1.02%  "

# Row 4
$ws.Range("E4").Value = "  +0.54%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.95"
$ws.Range("E5").Value = "  -1.02%  "

# Row 6
$ws.Range("E6").Value = "  -1.19%  "

# Row 7
$ws.Range("E7").Value = "  +0.54%  "

# Row 8
$ws.Range("E8").Value = "  -3.01%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0617"
$ws.Range("E9").Value = "  -3.53%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.14"
$ws.Range("E10").Value = "  -7.56%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0788"
$ws.Range("E11").Value = "  -0.93%  "

# Row 12
$ws.Range("D12").Value = "1.855.84"
$ws.Range("E12").Value = "  -0.81%  "

# Row 13
$ws.Range("D13").Value = "1.629.77"
$ws.Range("E13").Value = "  -2.53%  "

# Row 14
$ws.Range("E14").Value = "  -2.56%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.521"
$ws.Range("E15").Value = "  -4.23%  "

# Row 16
$ws.Range("D16").Value = "25.994.74"
$ws.Range("E16").Value = "  -0.44%  "

# Row 17
$ws.Range("D17").Value = "0.0₃0738"
$ws.Range("E17").Value = "  -3.49%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.17"
$ws.Range("E18").Value = "  -3.54%  "

# Row 19
$ws.Range("E19").Value = "  +0.48%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "189.46"
$ws.Range("E20").Value = "  -3.10%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.22"
$ws.Range("E21").Value = "  -3.18%  "

# Row 22
$ws.Range("E22").Value = "  -3.95%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.04"
$ws.Range("E23").Value = "  -2.90%  "

# Row 24
$ws.Range("E24").Value = "  +0.35%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.75"
$ws.Range("E25").Value = "  +0.30%  "

# Row 26
$ws.Range("B26").Value = "BinanceUSD"
$ws.Range("C26").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.01"
$ws.Range("E26").Value = "  +0.42%  "

# Row 27
$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.78"
$ws.Range("E27").Value = "  -1.14%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.70"
$ws.Range("E28").Value = "  -2.79%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.10"
$ws.Range("E29").Value = "  -2.93%  "

# Row 30
$ws.Range("E30").Value = "  -1.65%  "

# Row 31
$ws.Range("E31").Value = "  -3.48%  "

# Row 32
$ws.Range("E32").Value = "  -4.44%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.10"
$ws.Range("E33").Value = "  -5.83%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.41"
$ws.Range("E34").Value = "  -2.29%  "

# Row 35
$ws.Range("E35").Value = "  -3.36%  "

# Row 36
$ws.Range("D36").Value = "1.130.75"
$ws.Range("E36").Value = "  -0.21%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.848"
$ws.Range("E37").Value = "  -6.39%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.43"
$ws.Range("E38").Value = "  -1.42%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.514"
$ws.Range("E39").Value = "  -5.00%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0153"
$ws.Range("E40").Value = "  -2.15%  "

# Row 41
$ws.Range("E41").Value = "  -1.35%  "

# Row 42
$ws.Range("E42").Value = "  -3.21%  "

# Row 43
$ws.Range("D43").Value = "1.766.88"
$ws.Range("E43").Value = "  -0.78%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.18"
$ws.Range("E44").Value = "  -5.65%  "

# Row 45
$ws.Range("D45").Value = "0.0₆0112"
$ws.Range("E45").Value = "  -4.51%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "54.52"
$ws.Range("E46").Value = "  -3.91%  "

# Row 47
$ws.Range("E47").Value = "  +0.38%  "

# Row 48
$ws.Range("E48").Value = "  +0.24%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.414"
$ws.Range("E49").Value = "  -0.05%  "

# Row 50
$ws.Range("B50").Value = "USDD"
$ws.Range("C50").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.01"
$ws.Range("E50").Value = "  +0.55%  "

# Row 51
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.47"
$ws.Range("E51").Value = "  -4.06%  "
